# Weekly update: insert this week's new price record for
# "Hortaliza, Feria Lagunitas de Puerto Montt - Poroto verde" at the top of
# the data (row 17 of the historical series), pushing all older records
# down by one row (old row 17 becomes row 18, ..., old row 58 becomes the
# new row 59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 17..58 down to 18..59, leaving row 17 free for the new record.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with this week's record.
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C17").Value = "Los Lagos"
$ws.Range("D17").Value = 44607
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 100112031
$ws.Range("G17").Value = "Poroto verde"
$ws.Range("H17").Value = "Magnum"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 35000
$ws.Range("L17").Value = 35000
$ws.Range("M17").Value = 35000
$ws.Range("N17").Value = "$/saco 25 kilos"
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 1400
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
